$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" year column (Q) to the time series table, matching the
# formatting already used by the preceding year column (P).
$ws.Range("P4").Copy() | Out-Null
$ws.Range("Q4").PasteSpecial(-4122) | Out-Null
$ws.Range("Q4").Value = 2022

$ws.Range("P5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null
$ws.Range("Q5").Value = 64.2

$ws.Application.CutCopyMode = 0

# Update the active selection to match the post-edit workbook state.
$ws.Range("R4").Select()
